$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for all existing data rows (2-420)
#    from 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C420").Value = 45202

# 2. Row 420 gains an explicit row height (matches default 15, but becomes an explicit
#    customHeight entry in the saved XML).
$ws.Rows.Item(420).RowHeight = 15

# 3. Append the new record as row 421.
$ws.Range("A421").Value = "A 46326-2023"

$ws.Range("B421").Value = 45197
$ws.Range("B421").NumberFormat = "YYYY-MM-DD"

$ws.Range("C421").Value = 45202
$ws.Range("C421").NumberFormat = "YYYY-MM-DD"

$ws.Range("D421").Value = "VÄSTERBOTTENS LÄN"
$ws.Range("E421").Value = "NORDMALING"

$ws.Range("G421").Value = 4.6

$ws.Range("H421:Q421").Value = 0

$ws.Range("R421").Value = ""
$ws.Range("R421").WrapText = $true
